$wb = $excel.ActiveWorkbook

$statusNew = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6bfdbf6a5165793a6dae619fb53c9a037e694291/e2e/4f38d11f-f14e-400f-b6ba-7c0ecfb544fc.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/448c5631aebb469ff127846c0ff83e7d9922bb8b/e2e/4f38d11f-f14e-400f-b6ba-7c0ecfb544fc.md."

# --- Overview sheet: row 3 corresponds to 4f38d11f-f14e-400f-b6ba-7c0ecfb544fc.md ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusNew
$wsOverview.Range("F3").Value = $statusNew
$wsOverview.Range("G3").Value = "2016-08-15 14:46:20"

# --- zh-cn sheet: row 3 corresponds to 4f38d11f-f14e-400f-b6ba-7c0ecfb544fc.md ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusNew
$wsZhCn.Range("H3").Value = "2016-08-15 14:46:15"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

# --- de-de sheet: row 3 corresponds to 4f38d11f-f14e-400f-b6ba-7c0ecfb544fc.md ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusNew
$wsDeDe.Range("H3").Value = "2016-08-15 14:46:20"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
